$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new Hanuman / assistantExaminer row (row 29) --------------
# Seed it from row 28 so every column inherits that row's look, then
# overwrite with the real values/formula.
$ws.Range("A28:D28").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# B29 should carry B2's current (pre-fix) "outlier" look - grab it now,
# before B2 itself gets normalised below.
$ws.Range("B2").Copy()
$ws.Range("B29").PasteSpecial(-4122)

$ws.Range("A29").Value = "Hanuman"
$ws.Range("B29").Value = "0947241"
$ws.Range("C29").Value = "kurnool_eGov@123"
$ws.Range("D29").Formula = "=FALSE()"

$ws.Hyperlinks.Add($ws.Range("C29"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")

# --- Fix the mis-spelled / mis-prefixed dataName values in column A ---
$ws.Range("A23").Value = "Mohammed"
$ws.Range("A24").Value = "Ramachandra"

# --- B2 picks up the standard column-B formatting (same look as B3:B28) ---
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Value = "0944181"

# --- Update the view: scrolled down with the new row selected ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A30").Select()
